$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.594.10'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.07%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.841.86'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.16%  '

$ws.Range('E4').Value = '  +0.26%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.43'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.68%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.33%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4271'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.48%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3623'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.56%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07300'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.46%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8769'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.02%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.60'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.18%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.902.67'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.43%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.331'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.58%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.509'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.12%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06990'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.71%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.004'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.36%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '79.25'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.71%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008946'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.92%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.34'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.92%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '27.623.87'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.09%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.974'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.24%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.29'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.08%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.072.45'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.52%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.988'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.06%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '155.41'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.64%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.47'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.58%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '119.45'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.53%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.202'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.02%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.875'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.98%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08887'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.46%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7594'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.86%  '

$ws.Range('B33').Value = 'HuobiToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.948'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.20%  '

$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.501'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.41%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.124'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.63%  '

$ws.Range('E36').Value = '  +0.41%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.05450'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.89%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.108'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.64%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01932'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.11%  '

$ws.Range('E40').Value = '  +0.65%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1660'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.61%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5064'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.11%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.552'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.90%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.381'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.51%  '

$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.06552'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.13%  '

$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.39'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.32%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '105.78'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.63%  '

$ws.Range('B48').Value = 'PaxDollar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.001'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.40%  '

$ws.Range('B49').Value = 'Decentraland'
$ws.Range('C49').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4637'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.65%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.634'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.22%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '64.67'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.17%  '
